$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.604.69"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.876.11"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'247.63"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4748"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "'0.2906"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.06484"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'21.96"
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("D11").Value = "'0.07737"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "'0.7385"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "1.874.63"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "'95.98"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "'5.178"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'274.39"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").Value = "30.588.52"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'0.000007478"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "2.122.76"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'5.215"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'6.176"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'165.43"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'9.198"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "'18.78"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "'1.905"
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("D29").Value = "'0.09886"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "'4.251"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").Value = "'4.091"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'0.04779"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'1.121"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "'0.6942"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'0.01849"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").Value = "'2.759"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "'6.236"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").Value = "'73.24"
$ws.Range("E41").Value = "  +3.54%  "
$ws.Range("D42").Value = "'1.980"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").Value = "'0.4174"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.8358"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "'101.52"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").Value = "'9.352"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'35.38"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "'6.969"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").Value = "'918.91"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D51").Value = "'0.05666"
$ws.Range("E51").Value = "  +0.69%  "
